$d = $word.ActiveDocument

while ($d.Tables.Count -gt 0) {
    $d.Tables(1).Delete()
}

while ($d.Paragraphs.Count -gt 1) {
    $d.Paragraphs(2).Range.Delete()
}

$d.Paragraphs(1).Range.LanguageID = "en-IN"

$d.Save()
